$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)

$numRange = $hdr.Range.Duplicate
$numRange.SetRange(6,7)
Write-Host "numRange text before = [$($numRange.Text)]"

$openXmlTemplate = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="72 Light" w:hAnsi="72 Light" w:cs="72 Light"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

try {
    $numRange.WordOpenXML = $openXmlTemplate
    Write-Host "Set WordOpenXML worked"
} catch {
    Write-Host "Error: $_"
}
Write-Host "After: [$($hdr.Range.Text)]"
